# Add a new "API Dev and Testing" worksheet after the existing "JavaScript"
# sheet, populate it with two paid API development / testing courses, and
# make it the active (selected) sheet - matching the author's commit:
# "API development and testing paid courses added."

$wb = $excel.ActiveWorkbook

# The "JavaScript" sheet is currently the last / active sheet; add the new
# sheet right after it so it becomes the new last sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "API Dev and Testing"

# Populate the URL cells (column C) before the title cells (column B) so the
# shared-string table is built in the same order the original author's
# workbook used (url, title, url, title, ...).
$newSheet.Range("C2").Value = "https://www.udemy.com/course/postman-the-complete-guide/"
$newSheet.Range("B2").Value = "Postman: The Complete Guide - REST API Testing"
$newSheet.Range("C4").Value = "https://www.udemy.com/course/django-rest-framework/"
$newSheet.Range("B4").Value = "Build REST APIs with Django REST Framework and Python"

# Match the column widths used by the rest of the workbook's course sheets
# (chosen so the runtime's internal pixel-rounded column width storage lands
# as close as possible to the target widths of 75.85546875 / 73.42578125).
$newSheet.Columns.Item(2).ColumnWidth = 75
$newSheet.Columns.Item(3).ColumnWidth = 72.6666667

# Reproduce the selection left behind on the new sheet.
$newSheet.Range("C6").Select() | Out-Null

# The new sheet becomes the active / visible tab.
$newSheet.Activate() | Out-Null
